# Add three new characters (Akai, Yuka, SteampunkRobot) to the ActorTable
# sheet, rows 15-17 (Actor014/015/016), per commit:
# "신규 3종 캐릭터 테이블 추가 14 Akai 15 Yuka 16 SteampunkRobot"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ActorTable")

# --- prefabAddress (column L) filled first for all three new rows ---
$ws.Range("L15").Value = "Akai"
$ws.Range("L16").Value = "Yuka"
$ws.Range("L17").Value = "SteampunkRobot"

# --- Row 15: Actor014 / Akai ---
$ws.Range("B15").Value = "CharName_Akai"
$ws.Range("C15").Value = "CharDesc_Akai"
$ws.Range("M15").Value = "Portrait_Akai"
$ws.Range("D15").Value = 1
$ws.Range("I15").Value = 3

# --- Row 16: Actor015 / Yuka ---
$ws.Range("B16").Value = "CharName_Yuka"
$ws.Range("C16").Value = "CharDesc_Yuka"
$ws.Range("M16").Value = "Portrait_Yuka"
$ws.Range("D16").Value = 1
$ws.Range("I16").Value = 1

# --- Row 17: Actor016 / SteampunkRobot ---
$ws.Range("B17").Value = "CharName_SteampunkRobot"
$ws.Range("C17").Value = "CharDesc_SteampunkRobot"
$ws.Range("M17").Value = "Portrait_SteampunkRobot"
$ws.Range("I17").Value = 1

# --- Page setup (paper size / orientation now explicit on ActorTable) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
